$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.470.65'
$ws.Range("E2").Value = '  -0.65%  '
$ws.Range("D3").Value = '2.513.70'
$ws.Range("E3").Value = '  -1.87%  '
$ws.Range("E4").Value = '  +0.05%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '308.10'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +1.74%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '96.35'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -0.97%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.588'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +2.14%  '
$ws.Range("E8").Value = '  +0.00%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.538'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -1.58%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '36.64'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +0.47%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.0813'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +0.41%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '7.73'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +0.64%  '
$ws.Range("E13").Value = '  -3.62%  '
$ws.Range("D14").Value = '2.899.43'
$ws.Range("E14").Value = '  -1.91%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '15.77'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +8.43%  '
$ws.Range("D16").Value = '2.483.61'
$ws.Range("E16").Value = '  -3.50%  '
$ws.Range("E17").Value = '  -2.30%  '
$ws.Range("D18").Value = '42.463.60'
$ws.Range("E18").Value = '  -0.83%  '
$ws.Range("E19").Value = '  -4.75%  '
$ws.Range("E20").Value = '  -1.62%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '6.46'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -2.62%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '71.51'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -0.21%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '253.30'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -1.51%  '
$ws.Range("E24").Value = '  -1.06%  '
$ws.Range("E25").Value = '  -2.64%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '26.99'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -4.31%  '
$ws.Range("E27").Value = '  +0.03%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '2.32'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +10.06%  '
$ws.Range("E29").Value = '  +0.57%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '37.48'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -4.41%  '
$ws.Range("E31").Value = '  -1.01%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '153.88'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -1.69%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '19.14'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +5.27%  '
$ws.Range("E34").Value = '  -2.22%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.0787'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -1.85%  '
$ws.Range("E36").Value = '  -4.39%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '2.61'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -4.91%  '
$ws.Range("E38").Value = '  -1.11%  '
$ws.Range("E39").Value = '  +0.38%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '24.21'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -10.51%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '3.40'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +0.80%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '3.87'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +0.26%  '
$ws.Range("E43").Value = '  -0.21%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("E45").Value = '  -1.13%  '
$ws.Range("D46").Value = '2.036.41'
$ws.Range("E46").Value = '  -1.19%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '84.56'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -4.42%  '
$ws.Range("E48").Value = '  -3.41%  '
$ws.Range("D49").Value = '2.754.87'
$ws.Range("E49").Value = '  -2.04%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '73.00'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -5.20%  '
$ws.Range("E51").Value = '  -0.06%  '
